$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1: copy format (style) from an existing header cell, then set text
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# Data for columns I and J, rows 2-45
$iValues = @(7,5,8,8,8,9,7,6,7,9,6,6,8,5,8,7,8,7,6,6,6,6,6,12,10,9,5,8,6,6,7,5,7,6,5,7,8,7,8,6,7,8,3,5)
$jValues = @(7,6,8,8,8,9,8,7,7,9,7,7,8,6,8,7,9,8,6,6,6,6,6,12,11,9,5,8,6,6,7,6,8,7,5,7,8,7,8,6,7,8,3,5)

for ($r = 2; $r -le 45; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $iValues[$idx]
    $ws.Cells.Item($r, 10).Value = $jValues[$idx]
}
